# Apply the "better interpretability" update to the team-formation template.
#  1. Remove the now-obsolete "User Guide" sheet.
#  2. Expand the three rightmost header labels on "Data" with short legends
#     explaining the 1/2/3 coding scheme, and wrap the header row.
#  3. Widen the "Employment Status"/"Teamwork style" columns to fit the
#     longer headers and grow the header row height to match.
#  4. Fix six mis-coded "Teamwork style" values (were 1 or 3, should be 2).
#  5. Leave the cursor/selection the way the author left it when saving.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Data")

# --- 1. Drop the "User Guide" sheet -----------------------------------
$guide = $wb.Worksheets.Item("User Guide")
$guide.Delete()

# --- 2. Update the header row text with the legends --------------------
$ws.Range("D1").Value = "How familiar are you with Data Science?`n1 = Have not attended a DS class, 2 = Attended 1 or 2 DS classes, 3 = Attended 3+ DS classes"
$ws.Range("E1").Value = "Employment Status`n1 = Not employed, 2 = Part-time, 3 = Full-time professional"
$ws.Range("F1").Value = "Teamwork style`n1 = Careful, 2 = Clutch, 3 = Acceleration"

# Wrap the header text and grow the row so the legends are readable.
$ws.Range("A1:F1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 68

# --- 3. Resize the two affected columns ---------------------------------
$ws.Columns.Item(5).ColumnWidth = 25.83203125
$ws.Columns.Item(6).ColumnWidth = 28.5

# --- 4. Correct the six mis-coded "Teamwork style" values ---------------
$ws.Range("F10").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F30").Value = 2
$ws.Range("F33").Value = 2
$ws.Range("F39").Value = 2
$ws.Range("F46").Value = 2

# --- 5. Restore the author's on-save selection ---------------------------
$ws.Activate()
$ws.Range("E42").Select()
